# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Column G ("K") holds recalculated strikeout counts for each outing row (rows 2-66).
# Update every row whose recomputed K value differs from what is currently stored.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    3  = 0
    4  = 0
    5  = 2
    6  = 0
    7  = 0
    8  = 2
    10 = 1
    11 = 1
    12 = 1
    13 = 2
    14 = 1
    15 = 1
    16 = 2
    17 = 3
    18 = 1
    19 = 0
    20 = 0
    21 = 0
    22 = 0
    23 = 1
    24 = 2
    25 = 1
    26 = 2
    28 = 2
    29 = 1
    30 = 0
    31 = 2
    32 = 0
    33 = 0
    34 = 1
    35 = 1
    36 = 2
    37 = 1
    38 = 0
    39 = 3
    40 = 1
    41 = 0
    42 = 2
    43 = 1
    44 = 1
    45 = 2
    46 = 2
    47 = 1
    48 = 0
    49 = 1
    50 = 2
    51 = 2
    52 = 1
    53 = 3
    54 = 2
    55 = 2
    56 = 2
    57 = 2
    58 = 3
    59 = 2
    60 = 2
    61 = 3
    62 = 1
    63 = 1
    64 = 2
    65 = 2
    66 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
